# #5: cash & deposit done
# Rebuild the "存款" (deposit) worksheet (sheet4) so that row 1 holds the
# generic field-name headers (matching the other property sheets) and
# row 2 holds the fully-populated data row, including the new
# property_category / category / date / legislator_name / legislator_id /
# source_file / index columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Row 1: header labels -----------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Give the newly added header cells (G1:M1) the same bold/centered/bordered
# formatting already used by the existing header cells.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

# --- Row 2: data values ---------------------------------------------------
# B2:F2 already contain the correct values (bank name, deposit type,
# currency, owner, total amount) - only the new trailing columns need to be
# filled in.
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"

# "2012-04-30" must stay a plain text value (like every other date-ish
# column in this workbook) rather than being auto-converted into a date
# serial number. Build it via a text formula and paste only the resulting
# value so no new number-format style gets created.
$ws.Range("Z99").Formula = "=""2012-04-30"""
$ws.Range("Z99").Copy()
$ws.Range("I2").PasteSpecial(-4163)
$ws.Range("Z99").ClearContents()

$ws.Range("J2").Value = "何欣純"
$ws.Range("K2").Value = 1733
$ws.Range("L2").Value = "tmp2e891"
$ws.Range("M2").Value = 47

# Match the plain (unbordered/unbold) formatting of the rest of row 2.
$ws.Range("B2").Copy()
$ws.Range("G2:M2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
